$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four rows for the old "MuSCs" sending-cluster block (original rows 10-13),
# which no longer exist in the updated TPM output.
$ws.Range("A10:T13").EntireRow.Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna3"
$ws.Range("C2").Value = "Epha1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3193606666666667
$ws.Range("H2").Value = 0.958082
$ws.Range("I2").Value = 0.7979421849584948
$ws.Range("J2").Value = 0.7979421849584948
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.844648666666667
$ws.Range("N2").Value = 17.533946
$ws.Range("O2").Value = 0.3204643139023235
$ws.Range("P2").Value = 0.3204643139023235
$ws.Range("Q2").Value = 1.866550894619111
$ws.Range("R2").Value = 16.798958051572
$ws.Range("S2").Value = 0.255711994836445
$ws.Range("T2").Value = 0.255711994836445

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna3"
$ws.Range("C3").Value = "Epha1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3193606666666667
$ws.Range("H3").Value = 0.958082
$ws.Range("I3").Value = 0.7979421849584948
$ws.Range("J3").Value = 0.7979421849584948
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.370261666666667
$ws.Range("N3").Value = 13.110785
$ws.Range("O3").Value = 0.2396231127748354
$ws.Range("P3").Value = 0.2396231127748355
$ws.Range("Q3").Value = 1.395689679374445
$ws.Range("R3").Value = 12.56120711437
$ws.Range("S3").Value = 0.191205390174108
$ws.Range("T3").Value = 0.191205390174108

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna3"
$ws.Range("C4").Value = "Epha1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3193606666666667
$ws.Range("H4").Value = 0.958082
$ws.Range("I4").Value = 0.7979421849584948
$ws.Range("J4").Value = 0.7979421849584948
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.452372666666666
$ws.Range("N4").Value = 19.357118
$ws.Range("O4").Value = 0.3537860524377295
$ws.Range("P4").Value = 0.3537860524377295
$ws.Range("Q4").Value = 2.060634036408445
$ws.Range("R4").Value = 18.545706327676
$ws.Range("S4").Value = 0.2823008156900025
$ws.Range("T4").Value = 0.2823008156900025

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna3"
$ws.Range("C5").Value = "Epha1"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3193606666666667
$ws.Range("H5").Value = 0.958082
$ws.Range("I5").Value = 0.7979421849584948
$ws.Range("J5").Value = 0.7979421849584948
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.570781
$ws.Range("N5").Value = 4.712343
$ws.Range("O5").Value = 0.08612652088511148
$ws.Range("P5").Value = 0.0861265208851115
$ws.Range("Q5").Value = 0.5016456673473334
$ws.Range("R5").Value = 4.514811006126
$ws.Range("S5").Value = 0.0687239842579393
$ws.Range("T5").Value = 0.0687239842579393

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efna3"
$ws.Range("C6").Value = "Epha1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.08086966666666666
$ws.Range("H6").Value = 0.242609
$ws.Range("I6").Value = 0.2020578150415052
$ws.Range("J6").Value = 0.2020578150415053
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.844648666666667
$ws.Range("N6").Value = 17.533946
$ws.Range("O6").Value = 0.3204643139023235
$ws.Range("P6").Value = 0.3204643139023235
$ws.Range("Q6").Value = 0.4726547894571111
$ws.Range("R6").Value = 4.253893105114
$ws.Range("S6").Value = 0.06475231906587857
$ws.Range("T6").Value = 0.06475231906587858

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna3"
$ws.Range("C7").Value = "Epha1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.08086966666666666
$ws.Range("H7").Value = 0.242609
$ws.Range("I7").Value = 0.2020578150415052
$ws.Range("J7").Value = 0.2020578150415053
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.370261666666667
$ws.Range("N7").Value = 13.110785
$ws.Range("O7").Value = 0.2396231127748354
$ws.Range("P7").Value = 0.2396231127748355
$ws.Range("Q7").Value = 0.3534216042294444
$ws.Range("R7").Value = 3.180794438065
$ws.Range("S7").Value = 0.04841772260072745
$ws.Range("T7").Value = 0.04841772260072746

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna3"
$ws.Range("C8").Value = "Epha1"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.08086966666666666
$ws.Range("H8").Value = 0.242609
$ws.Range("I8").Value = 0.2020578150415052
$ws.Range("J8").Value = 0.2020578150415053
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.452372666666666
$ws.Range("N8").Value = 19.357118
$ws.Range("O8").Value = 0.3537860524377295
$ws.Range("P8").Value = 0.3537860524377295
$ws.Range("Q8").Value = 0.5218012267624443
$ws.Range("R8").Value = 4.696211040862
$ws.Range("S8").Value = 0.07148523674772701
$ws.Range("T8").Value = 0.07148523674772704

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna3"
$ws.Range("C9").Value = "Epha1"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.08086966666666666
$ws.Range("H9").Value = 0.242609
$ws.Range("I9").Value = 0.2020578150415052
$ws.Range("J9").Value = 0.2020578150415053
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.570781
$ws.Range("N9").Value = 4.712343
$ws.Range("O9").Value = 0.08612652088511148
$ws.Range("P9").Value = 0.0861265208851115
$ws.Range("Q9").Value = 0.1270285358763333
$ws.Range("R9").Value = 1.143256822887
$ws.Range("S9").Value = 0.0174025366271722
$ws.Range("T9").Value = 0.0174025366271722

Write-Output "applied"